# Auto-generated edit script: refresh market-price derived columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR profit tables, per scheduled-runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 44776.2
$ws.Range("I6").Value = 86.5
$ws.Range("J6").Value = 111810.75
$ws.Range("K6").Value = 259.5
$ws.Range("L6").Value = 335432.25
$ws.Range("M6").Value = -147.5
$ws.Range("N6").Value = -335656.25
$ws.Range("H28").Value = 888.8333
$ws.Range("I28").Value = 1025.5555
$ws.Range("J28").Value = 752.1111
$ws.Range("K28").Value = 1025.5555
$ws.Range("L28").Value = 752.1111
$ws.Range("M28").Value = -540.5554999999999
$ws.Range("N28").Value = -1722.1111
$ws.Range("H131").Value = 1281.8
$ws.Range("I131").Value = 964.0714
$ws.Range("J131").Value = 1493.619
$ws.Range("K131").Value = 2892.2142
$ws.Range("L131").Value = 4480.857
$ws.Range("M131").Value = 2147.7858
$ws.Range("N131").Value = -14560.857
$ws.Range("H138").Value = 2944159.5
$ws.Range("I138").Value = 2226.2646
$ws.Range("J138").Value = 5886092.5
$ws.Range("K138").Value = 6678.793799999999
$ws.Range("L138").Value = 17658277.5
$ws.Range("M138").Value = -1538.793799999999
$ws.Range("N138").Value = -17668557.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 40083060
$ws.Range("I61").Value = 47668500
$ws.Range("J61").Value = 259500
$ws.Range("K61").Value = 47668500
$ws.Range("L61").Value = 259500
$ws.Range("M61").Value = -47668288
$ws.Range("N61").Value = -259924
$ws.Range("H109").Value = 34350
$ws.Range("J109").Value = 34350
$ws.Range("L109").Value = 34350
$ws.Range("N109").Value = -37124
$ws.Range("H110").Value = 1526.625
$ws.Range("I110").Value = 537.4
$ws.Range("J110").Value = 3175.3333
$ws.Range("K110").Value = 537.4
$ws.Range("L110").Value = 3175.3333
$ws.Range("M110").Value = 1507.6
$ws.Range("N110").Value = -7265.3333
$ws.Range("H122").Value = 4631530
$ws.Range("I122").Value = 1845.6818
$ws.Range("K122").Value = 5537.0454
$ws.Range("M122").Value = -3087.0454
$ws.Range("H132").Value = 17934324
$ws.Range("I132").Value = 27835850
$ws.Range("J132").Value = 111572.8
$ws.Range("K132").Value = 83507550
$ws.Range("L132").Value = 334718.4
$ws.Range("M132").Value = -83505020
$ws.Range("N132").Value = -339778.4
$ws.Range("H136").Value = 40083060
$ws.Range("I136").Value = 47668500
$ws.Range("J136").Value = 259500
$ws.Range("K136").Value = 143005500
$ws.Range("L136").Value = 778500
$ws.Range("M136").Value = -143002950
$ws.Range("N136").Value = -783600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 38463976
$ws.Range("I105").Value = 62501460
$ws.Range("J105").Value = 3998
$ws.Range("K105").Value = 62501460
$ws.Range("L105").Value = 3998
$ws.Range("M105").Value = -62499713
$ws.Range("N105").Value = -7492

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1797.0769
$ws.Range("I16").Value = 1832.2
$ws.Range("K16").Value = 1832.2
$ws.Range("M16").Value = -1545.2
$ws.Range("H31").Value = 62494.8
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 62494.8
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 62494.8
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -63084.8
$ws.Range("H34").Value = 62494.8
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 62494.8
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 62494.8
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -62898.8
$ws.Range("H113").Value = 1797.0769
$ws.Range("I113").Value = 1832.2
$ws.Range("K113").Value = 1832.2
$ws.Range("M113").Value = 337.8
$ws.Range("H134").Value = 49924
$ws.Range("I134").Value = 789.86664
$ws.Range("J134").Value = 155211.42
$ws.Range("K134").Value = 2369.59992
$ws.Range("L134").Value = 465634.26
$ws.Range("M134").Value = 165.4000800000003
$ws.Range("N134").Value = -470704.26

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1015
$ws.Range("J131").Value = 1050
$ws.Range("L131").Value = 3150
$ws.Range("N131").Value = -13230
$ws.Range("H136").Value = 3164.375
$ws.Range("I136").Value = 2279.8
$ws.Range("J136").Value = 4638.6665
$ws.Range("K136").Value = 6839.400000000001
$ws.Range("L136").Value = 13915.9995
$ws.Range("M136").Value = -1739.400000000001
$ws.Range("N136").Value = -24115.9995
$ws.Range("H138").Value = 5091.7646
$ws.Range("I138").Value = 1278.1818
$ws.Range("J138").Value = 12083.333
$ws.Range("K138").Value = 3834.5454
$ws.Range("L138").Value = 36249.999
$ws.Range("M138").Value = 1305.4546
$ws.Range("N138").Value = -46529.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1850.0714
$ws.Range("I61").Value = 2019.4667
$ws.Range("J61").Value = 1654.6154
$ws.Range("K61").Value = 2019.4667
$ws.Range("L61").Value = 1654.6154
$ws.Range("M61").Value = -1817.4667
$ws.Range("N61").Value = -2058.6154
$ws.Range("H68").Value = 1798.8
$ws.Range("I68").Value = 1704.1177
$ws.Range("J68").Value = 2335.3333
$ws.Range("K68").Value = 1704.1177
$ws.Range("L68").Value = 2335.3333
$ws.Range("M68").Value = -955.1177
$ws.Range("N68").Value = -3833.3333
$ws.Range("H71").Value = 1798.8
$ws.Range("I71").Value = 1704.1177
$ws.Range("J71").Value = 2335.3333
$ws.Range("K71").Value = 8520.5885
$ws.Range("L71").Value = 11676.6665
$ws.Range("M71").Value = -4776.5885
$ws.Range("N71").Value = -19164.6665
$ws.Range("H82").Value = 1999
$ws.Range("I82").Value = 960.8333
$ws.Range("K82").Value = 960.8333
$ws.Range("M82").Value = -599.8333
$ws.Range("H85").Value = 1999
$ws.Range("I85").Value = 960.8333
$ws.Range("K85").Value = 960.8333
$ws.Range("M85").Value = 287.1667
$ws.Range("H93").Value = 1617
$ws.Range("I93").Value = 1574.4445
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 1574.4445
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = -326.4445000000001
$ws.Range("N93").Value = -4496
$ws.Range("H100").Value = 1894.2354
$ws.Range("I100").Value = 1633.6666
$ws.Range("K100").Value = 1633.6666
$ws.Range("M100").Value = -1092.6666
$ws.Range("H113").Value = 1850.0714
$ws.Range("I113").Value = 2019.4667
$ws.Range("J113").Value = 1654.6154
$ws.Range("K113").Value = 2019.4667
$ws.Range("L113").Value = 1654.6154
$ws.Range("M113").Value = 150.5333000000001
$ws.Range("N113").Value = -5994.6154
$ws.Range("H132").Value = 20822.574
$ws.Range("I132").Value = 2160.1177
$ws.Range("K132").Value = 6480.353099999999
$ws.Range("M132").Value = -3950.353099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 447.45715
$ws.Range("I107").Value = 359.25
$ws.Range("K107").Value = 1077.75
$ws.Range("M107").Value = 842.25
$ws.Range("H138").Value = 48012.43
$ws.Range("J138").Value = 48012.43
$ws.Range("L138").Value = 48012.43
$ws.Range("N138").Value = -58292.43
